$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 no longer exists in the updated report - remove it
$ws.Rows.Item(22).Delete()

# Row 2
$ws.Range("A2").Value = [double]"1.9623672415036708E-4"
$ws.Range("B2").Value = "firConvolutionLoopPipelining_top"

# Row 3
$ws.Range("A3").Value = [double]"2.8202610337757505E-5"
$ws.Range("B3").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_30"
$ws.Range("C3").Value = [double]"33.043479919433594"
$ws.Range("D3").Value = [double]"5.609173774719238"
$ws.Range("E3").Value = [double]"32.0"
$ws.Range("F3").Value = [double]"9.0"
$ws.Range("G3").Value = "myclk"
$ws.Range("H3").Value = "FF "

# Row 4
$ws.Range("A4").Value = [double]"2.6129782781936228E-5"
$ws.Range("B4").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_00"
$ws.Range("C4").Value = [double]"33.4782600402832"
$ws.Range("D4").Value = [double]"5.349087238311768"
$ws.Range("E4").Value = [double]"32.0"
$ws.Range("F4").Value = [double]"11.0"
$ws.Range("G4").Value = "myclk"
$ws.Range("H4").Value = "FF "

# Row 5
$ws.Range("A5").Value = [double]"1.8403477952233516E-5"
$ws.Range("B5").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_90"
$ws.Range("C5").Value = [double]"11.30434799194336"
$ws.Range("D5").Value = [double]"5.65217399597168"
$ws.Range("E5").Value = [double]"32.0"
$ws.Range("F5").Value = [double]"20.0"
$ws.Range("G5").Value = "myclk"
$ws.Range("H5").Value = "FF "

# Row 6
$ws.Range("A6").Value = [double]"1.71486954059219E-5"
$ws.Range("B6").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_80"
$ws.Range("C6").Value = [double]"22.60869598388672"
$ws.Range("D6").Value = [double]"5.65217399597168"
$ws.Range("E6").Value = [double]"32.0"
$ws.Range("F6").Value = [double]"8.0"
$ws.Range("G6").Value = "myclk"
$ws.Range("H6").Value = "FF "

# Row 7
$ws.Range("A7").Value = [double]"1.3938260963186622E-5"
$ws.Range("B7").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_70"
$ws.Range("C7").Value = [double]"22.60869598388672"
$ws.Range("D7").Value = [double]"5.65217399597168"
$ws.Range("E7").Value = [double]"32.0"
$ws.Range("F7").Value = [double]"7.0"
$ws.Range("G7").Value = "myclk"
$ws.Range("H7").Value = "FF "

# Row 8
$ws.Range("A8").Value = [double]"1.367826098430669E-5"
$ws.Range("B8").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_40"
$ws.Range("C8").Value = [double]"22.60869598388672"
$ws.Range("D8").Value = [double]"5.65217399597168"
$ws.Range("E8").Value = [double]"32.0"
$ws.Range("F8").Value = [double]"8.0"
$ws.Range("G8").Value = "myclk"
$ws.Range("H8").Value = "FF "

# Row 9
$ws.Range("A9").Value = [double]"1.3093695997667965E-5"
$ws.Range("B9").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_10"
$ws.Range("C9").Value = [double]"22.173913955688477"
$ws.Range("D9").Value = [double]"5.503129959106445"
$ws.Range("E9").Value = [double]"32.0"
$ws.Range("F9").Value = [double]"6.0"
$ws.Range("G9").Value = "myclk"
$ws.Range("H9").Value = "FF "

# Row 10
$ws.Range("A10").Value = [double]"1.0422609193483368E-5"
$ws.Range("B10").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_50"
$ws.Range("C10").Value = [double]"11.30434799194336"
$ws.Range("D10").Value = [double]"5.65217399597168"
$ws.Range("E10").Value = [double]"32.0"
$ws.Range("F10").Value = [double]"12.0"
$ws.Range("G10").Value = "myclk"
$ws.Range("H10").Value = "FF "

# Row 11
$ws.Range("A11").Value = [double]"7.986522177816369E-6"
$ws.Range("B11").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_20"
$ws.Range("C11").Value = [double]"11.30434799194336"
$ws.Range("D11").Value = [double]"5.65217399597168"
$ws.Range("E11").Value = [double]"32.0"
$ws.Range("F11").Value = [double]"8.0"
$ws.Range("G11").Value = "myclk"
$ws.Range("H11").Value = "FF "

# Row 12
$ws.Range("A12").Value = [double]"7.92434821050847E-6"
$ws.Range("B12").Value = "firConvolutionLoopPipelining_IP/U0/shiftRegister_60"
$ws.Range("C12").Value = [double]"11.30434799194336"
$ws.Range("D12").Value = [double]"5.65217399597168"
$ws.Range("E12").Value = [double]"32.0"
$ws.Range("F12").Value = [double]"7.0"
$ws.Range("G12").Value = "myclk"
$ws.Range("H12").Value = "FF "

# Row 13
$ws.Range("A13").Value = [double]"7.206521786429221E-6"
$ws.Range("B13").Value = "firConvolutionLoopPipelining_IP/U0/accumulator_reg_1140"
$ws.Range("C13").Value = [double]"10.8695650100708"
$ws.Range("D13").Value = [double]"60.71390914916992"
$ws.Range("E13").Value = [double]"32.0"
$ws.Range("F13").Value = [double]"8.0"
$ws.Range("G13").Value = "myclk"
$ws.Range("H13").Value = "FF "

# Row 14
$ws.Range("A14").Value = [double]"6.816521818109322E-6"
$ws.Range("B14").Value = "firConvolutionLoopPipelining_IP/U0/ap_NS_fsm1"
$ws.Range("C14").Value = [double]"11.30434799194336"
$ws.Range("D14").Value = [double]"5.533477783203125"
$ws.Range("E14").Value = [double]"32.0"
$ws.Range("F14").Value = [double]"10.0"
$ws.Range("G14").Value = "myclk"
$ws.Range("H14").Value = "FF "

# Row 15
$ws.Range("A15").Value = [double]"6.6358693402435165E-6"
$ws.Range("B15").Value = "firConvolutionLoopPipelining_IP/U0/ap_enable_reg_pp0_iter2"
$ws.Range("C15").Value = [double]"10.8695650100708"
$ws.Range("D15").Value = [double]"61.279170989990234"
$ws.Range("E15").Value = [double]"36.0"
$ws.Range("F15").Value = [double]"10.0"
$ws.Range("G15").Value = "myclk"
$ws.Range("H15").Value = "FF LUT "

# Row 16
$ws.Range("A16").Value = [double]"6.203478278621333E-6"
$ws.Range("B16").Value = "firConvolutionLoopPipelining_IP/U0/tmp_6_reg_463[16]_i_1_n_2"
$ws.Range("C16").Value = [double]"10.434782981872559"
$ws.Range("D16").Value = [double]"68.64195251464844"
$ws.Range("E16").Value = [double]"18.0"
$ws.Range("F16").Value = [double]"6.0"
$ws.Range("G16").Value = "myclk"
$ws.Range("H16").Value = "DSP FF "

# Row 17
$ws.Range("A17").Value = [double]"5.309782409312902E-6"
$ws.Range("B17").Value = "firConvolutionLoopPipelining_IP/U0/ap_enable_reg_pp0_iter1"
$ws.Range("C17").Value = [double]"10.8695650100708"
$ws.Range("D17").Value = [double]"61.713958740234375"
$ws.Range("E17").Value = [double]"12.0"
$ws.Range("F17").Value = [double]"7.0"
$ws.Range("G17").Value = "myclk"
$ws.Range("H17").Value = "FF LUT "

# Row 18
$ws.Range("A18").Value = [double]"3.3532608085806714E-6"
$ws.Range("B18").Value = "firConvolutionLoopPipelining_IP/U0/coefficientsFilter1_1_reg_4580"
$ws.Range("C18").Value = [double]"10.8695650100708"
$ws.Range("D18").Value = [double]"56.80443572998047"
$ws.Range("E18").Value = [double]"2.0"
$ws.Range("F18").Value = [double]"4.0"
$ws.Range("G18").Value = "myclk"
$ws.Range("H18").Value = "DSP "

# Row 19
$ws.Range("A19").Value = [double]"2.0815216430492E-6"
$ws.Range("B19").Value = "firConvolutionLoopPipelining_IP/U0/ce0"
$ws.Range("C19").Value = [double]"10.8695650100708"
$ws.Range("D19").Value = [double]"67.27208709716797"
$ws.Range("E19").Value = [double]"2.0"
$ws.Range("F19").Value = [double]"2.0"
$ws.Range("G19").Value = "myclk"
$ws.Range("H19").Value = "DSP "

# Row 20
$ws.Range("A20").Value = [double]"1.1249999261053745E-6"
$ws.Range("B20").Value = "firConvolutionLoopPipelining_IP/U0/i_reg_1270"
$ws.Range("C20").Value = [double]"10.8695650100708"
$ws.Range("D20").Value = [double]"62.10765075683594"
$ws.Range("E20").Value = [double]"5.0"
$ws.Range("F20").Value = [double]"1.0"
$ws.Range("G20").Value = "myclk"
$ws.Range("H20").Value = "FF "

# Row 21
$ws.Range("A21").Value = [double]"5.765217565567582E-7"
$ws.Range("B21").Value = "firConvolutionLoopPipelining_IP/U0/tmp_4_reg_4300"
$ws.Range("C21").Value = [double]"11.30434799194336"
$ws.Range("D21").Value = [double]"55.85365295410156"
$ws.Range("E21").Value = [double]"4.0"
$ws.Range("F21").Value = [double]"1.0"
$ws.Range("G21").Value = "myclk"
$ws.Range("H21").Value = "FF "
